# Weekly crime-data refresh for the 1st Precinct CompStat report.
# Updates the report header (volume/week numbers) and all weekly crime
# statistics rows (14-30) to the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: bump the report "Number" and the covered week dates.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 40

# ---------------------------------------------------------------------
# Row 16 - Robbery (D16 & E16 switch from numbers to the "0"/"***.*"
# placeholder text used elsewhere in the sheet for not-applicable data)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "***.*"
$ws.Range("N22").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 132
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -8.965517241379
$ws.Range("M16").Value = 41.935483870967
$ws.Range("N16").Value = -84.210526315789

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 83.333333333333
$ws.Range("I17").Value = 128
$ws.Range("J17").Value = 116
$ws.Range("K17").Value = 10.344827586206
$ws.Range("L17").Value = 58.024691358024
$ws.Range("M17").Value = 88.235294117647
$ws.Range("N17").Value = -27.272727272727

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 246
$ws.Range("J18").Value = 180
$ws.Range("K18").Value = 36.666666666666
$ws.Range("L18").Value = 3.797468354430
$ws.Range("M18").Value = 47.305389221556
$ws.Range("N18").Value = -71.853546910755

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 24
$ws.Range("E19").Value = 26.315789473684
$ws.Range("F19").Value = 90
$ws.Range("G19").Value = 90
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 1184
$ws.Range("J19").Value = 819
$ws.Range("K19").Value = 44.566544566544
$ws.Range("L19").Value = 82.998454404945
$ws.Range("M19").Value = 14.285714285714
$ws.Range("N19").Value = -68.552456839309

# ---------------------------------------------------------------------
# Row 20 - G.L.A. (C20 switches from a number to the "0" placeholder
# text used elsewhere in the sheet for not-applicable data)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 100
$ws.Range("J20").Value = 68
$ws.Range("K20").Value = 4.411764705882
$ws.Range("M20").Value = 86.842105263157
$ws.Range("N20").Value = -91.857798165137

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 38.461538461538
$ws.Range("F21").Value = 132
$ws.Range("G21").Value = 140
$ws.Range("H21").Value = -5.714285714285
$ws.Range("I21").Value = 1782
$ws.Range("J21").Value = 1332
$ws.Range("K21").Value = 33.783783783783
$ws.Range("L21").Value = 53.488372093023
$ws.Range("M21").Value = 26.5625
$ws.Range("N21").Value = -72.773109243697

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("I22").Value = 93
$ws.Range("J22").Value = 79
$ws.Range("K22").Value = 17.721518987341
$ws.Range("L22").Value = 10.714285714285
$ws.Range("M22").Value = 36.764705882352

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 88
$ws.Range("D24").Value = 65
$ws.Range("E24").Value = 35.384615384615
$ws.Range("F24").Value = 344
$ws.Range("G24").Value = 240
$ws.Range("H24").Value = 43.333333333333
$ws.Range("I24").Value = 3901
$ws.Range("J24").Value = 2196
$ws.Range("K24").Value = 77.641165755919
$ws.Range("L24").Value = 134.858518964479
$ws.Range("M24").Value = 137.720901889092

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = -5.263157894736
$ws.Range("I25").Value = 368
$ws.Range("J25").Value = 356
$ws.Range("K25").Value = 3.370786516853
$ws.Range("L25").Value = 47.791164658634
$ws.Range("M25").Value = 42.084942084942

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = 47.368421052631

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 113
$ws.Range("J27").Value = 82
$ws.Range("K27").Value = 37.804878048780
$ws.Range("L27").Value = 71.212121212121

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes (D30 & E30 switch from the "0"/"***.*" placeholder
# text to actual numbers now that data is available)
# ---------------------------------------------------------------------
$ws.Range("D30").Value = 4
$ws.Range("D20").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("E20").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 7
$ws.Range("J30").Value = 18
$ws.Range("K30").Value = -16.666666666666
